# Edit: add "Used"/"Location" tracking columns to the images sheet,
# mark the favicon row, widen/wrap description text, and make the
# images sheet the active tab.

$wb = $excel.ActiveWorkbook
$wsContent = $wb.Worksheets.Item("Content")
$wsImages  = $wb.Worksheets.Item("images")

# H9 ("icon_" row) description changes from "icon" to "favicon"
$wsImages.Range("H9").Value = "favicon"

# Mark usage for the vector volleyball row and the favicon row
$wsImages.Range("I8").Value = "Yes"
$wsImages.Range("I9").Value = "Yes"

# --- images sheet: new columns I (Used) and J (Location) ---
$wsImages.Range("I1").Value = "Used"
$wsImages.Range("J1").Value = "Location"

$wsImages.Range("J9").Value = "All Pages"
$wsImages.Range("J8").Value = "Logo"

# --- column widths / wrapping ---
$wsImages.Columns.Item(2).ColumnWidth = 23.85546875
$wsImages.Columns.Item(10).ColumnWidth = 11.28515625

$wsImages.Range("A1:J9").WrapText = $true

# Row heights grow to fit the now-wrapped, narrower description column
$wsImages.Rows.Item(1).RowHeight = 37.5
$wsImages.Rows.Item(2).RowHeight = 75
$wsImages.Rows.Item(3).RowHeight = 37.5
$wsImages.Rows.Item(8).RowHeight = 37.5

# --- make images the active sheet/tab ---
$wsImages.Activate()
$wsContent.Range("E16").Select() | Out-Null
$wsImages.Range("D13").Select() | Out-Null

# --- page setup for images sheet ---
$wsImages.PageSetup.PaperSize = 9
$wsImages.PageSetup.Orientation = 1
